# Insert a new data row at sheet row 107 (pushes existing rows 107-205 down to 108-206)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(107).Insert()

# Populate the newly inserted row 107 with the new record's data
$ws.Cells.Item(107, 1).Value = 10
$ws.Cells.Item(107, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(107, 3).Value = "La Araucanía"
$ws.Cells.Item(107, 4).Value = 44566
$ws.Cells.Item(107, 5).Value = 9
$ws.Cells.Item(107, 6).Value = "Fruta"
$ws.Cells.Item(107, 7).Value = 100103
$ws.Cells.Item(107, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(107, 9).Value = 100103004
$ws.Cells.Item(107, 10).Value = "Durazno"
$ws.Cells.Item(107, 11).Value = "Early Majestic"
$ws.Cells.Item(107, 12).Value = "Especial"
$ws.Cells.Item(107, 13).Value = 155
$ws.Cells.Item(107, 14).Value = 25000
$ws.Cells.Item(107, 15).Value = 25000
$ws.Cells.Item(107, 16).Value = 25000
$ws.Cells.Item(107, 17).Value = "$/caja 20 kilos empedrada"
$ws.Cells.Item(107, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(107, 19).Value = 1250
$ws.Cells.Item(107, 20).Value = 20
